$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "nome_cartao" -> "nome_cartao_credito"
$ws.Range("B1").Value = "nome_cartao_credito"

# Update the sample data row to reflect new card name and value
$ws.Range("E2").Value = "Supermercados ABC"
$ws.Range("H2").Value = 400
